# Apply updated cryptocurrency price/volume figures scraped on
# Wed Jul 31 04:33:41 UTC 2024 (GitHub Actions cron update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.825.07'
$ws.Range("E2").Value = '  -1.37%  '

$ws.Range("D3").Value = '''3.280.85'
$ws.Range("E3").Value = '  -1.11%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''584.46'
$ws.Range("E5").Value = '  +1.94%  '

$ws.Range("D6").Value = '''179.87'
$ws.Range("E6").Value = '  -1.68%  '

$ws.Range("D7").Value = '''0.652'
$ws.Range("E7").Value = '  +8.34%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  -3.73%  '

$ws.Range("D10").Value = '''6.75'
$ws.Range("E10").Value = '  +1.81%  '

$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("D12").Value = '''3.847.20'
$ws.Range("E12").Value = '  -1.22%  '

$ws.Range("E13").Value = '  -4.48%  '

$ws.Range("D14").Value = '''65.923.80'
$ws.Range("E14").Value = '  -1.26%  '

$ws.Range("D15").Value = '''26.16'
$ws.Range("E15").Value = '  -3.95%  '

$ws.Range("E16").Value = '  -2.61%  '

$ws.Range("D17").Value = '''3.234.99'
$ws.Range("E17").Value = '  -2.02%  '

$ws.Range("D18").Value = '''429.55'
$ws.Range("E18").Value = '  -1.83%  '

$ws.Range("D19").Value = '''13.22'
$ws.Range("E19").Value = '  -4.06%  '

$ws.Range("D20").Value = '''5.48'
$ws.Range("E20").Value = '  -3.63%  '

$ws.Range("D21").Value = '''7.38'
$ws.Range("E21").Value = '  -3.45%  '

$ws.Range("D22").Value = '''72.05'
$ws.Range("E22").Value = '  -2.58%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("D25").Value = '''3.433.30'
$ws.Range("E25").Value = '  -0.83%  '

$ws.Range("E26").Value = '  -1.29%  '

$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '''0.196'
$ws.Range("E27").Value = '  +1.12%  '

$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '''0.0000113'
$ws.Range("E28").Value = '  -4.93%  '

$ws.Range("D29").Value = '''8.90'
$ws.Range("E29").Value = '  -2.12%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("E31").Value = '  +0.14%  '

$ws.Range("D32").Value = '''22.28'
$ws.Range("E32").Value = '  -2.78%  '

$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("D34").Value = '''5.16'
$ws.Range("E34").Value = '  -3.33%  '

$ws.Range("D35").Value = '''6.59'
$ws.Range("E35").Value = '  -2.94%  '

$ws.Range("E36").Value = '  -3.55%  '

$ws.Range("E37").Value = '  -1.09%  '

$ws.Range("E38").Value = '  -6.14%  '

$ws.Range("D39").Value = '''26.37'
$ws.Range("E39").Value = '  -3.73%  '

$ws.Range("D40").Value = '''1.79'
$ws.Range("E40").Value = '  -3.96%  '

$ws.Range("D41").Value = '''2.781.41'
$ws.Range("E41").Value = '  -1.02%  '

$ws.Range("E42").Value = '  -3.13%  '

$ws.Range("D43").Value = '''4.32'
$ws.Range("E43").Value = '  -3.16%  '

$ws.Range("E44").Value = '  -0.32%  '

$ws.Range("D45").Value = '''0.0658'
$ws.Range("E45").Value = '  -2.89%  '

$ws.Range("E46").Value = '  -5.29%  '

$ws.Range("D47").Value = '''2.29'
$ws.Range("E47").Value = '  -2.45%  '

$ws.Range("D48").Value = '''314.97'
$ws.Range("E48").Value = '  -1.26%  '

$ws.Range("D49").Value = '''23.10'
$ws.Range("E49").Value = '  -5.09%  '

$ws.Range("E50").Value = '  -2.29%  '

$ws.Range("E51").Value = '  +6.08%  '

# The leading apostrophe above marks the cell as text (quote-prefixed),
# which also tags it with a "quote prefix" cell style. Reset the style
# back to Normal so only the cell VALUES change, keeping formatting
# identical to the original workbook.
$ws.Range("D2:D51").Style = "Normal"
